# Item.xlsx edit: add a new "Equip_Weapon_1" (开山斧) row to the Item table,
# and nudge the saved-selection cursor (matches the authored XML diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 9) -------------------------------------------------
$ws.Cells.Item(9, 1).Value  = "Equip_Weapon_1"   # ID
$ws.Cells.Item(9, 2).Value  = 1                   # ItemType
$ws.Cells.Item(9, 3).Value  = 1                   # ItemSubType
$ws.Cells.Item(9, 4).Value  = 1                   # Level
$ws.Cells.Item(9, 5).Value  = "开山斧"            # ShowName
$ws.Cells.Item(9, 6).Value  = "开山斧武器"        # Desc

# Icon is a numeric-looking code that must stay text (same as the other
# rows' Icon column, which is formatted as Text via NumberFormat "@").
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value  = "50004"             # Icon

$ws.Cells.Item(9, 8).Value  = 0                   # CoolDownTime
$ws.Cells.Item(9, 9).Value  = 10000               # OverlayCount
$ws.Cells.Item(9, 10).Value = 100                 # BuyPrice
$ws.Cells.Item(9, 11).Value = 100                 # SalePrice

# --- Grow the XML-mapped table so it covers the new row --------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K9"))

# --- Move the saved selection cursor (as recorded in the sheet view) -------
$ws.Range("K13").Select() | Out-Null
